$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 736.2222
$ws.Range("I4").Value = 440.5
$ws.Range("J4").Value = 1327.6666
$ws.Range("K4").Value = 440.5
$ws.Range("L4").Value = 1327.6666
$ws.Range("M4").Value = -326.5
$ws.Range("N4").Value = -1555.6666
$ws.Range("H33").Value = 201.45454
$ws.Range("I33").Value = 201.45454
$ws.Range("K33").Value = 201.45454
$ws.Range("M33").Value = 27.54545999999999
$ws.Range("H41").Value = 109
$ws.Range("I41").Value = 34.4
$ws.Range("J41").Value = 233.33333
$ws.Range("K41").Value = 34.4
$ws.Range("L41").Value = 233.33333
$ws.Range("M41").Value = 405.6
$ws.Range("N41").Value = -1113.33333
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H76").Value = 4755.121
$ws.Range("I76").Value = 3718.2856
$ws.Range("K76").Value = 3718.2856
$ws.Range("M76").Value = -3403.2856
$ws.Range("H79").Value = 4755.121
$ws.Range("I79").Value = 3718.2856
$ws.Range("K79").Value = 3718.2856
$ws.Range("M79").Value = -2626.2856
$ws.Range("H88").Value = 1894.3334
$ws.Range("I88").Value = 2835
$ws.Range("J88").Value = 1222.4286
$ws.Range("K88").Value = 2835
$ws.Range("L88").Value = 1222.4286
$ws.Range("M88").Value = -2429
$ws.Range("N88").Value = -2034.4286
$ws.Range("H91").Value = 1894.3334
$ws.Range("I91").Value = 2835
$ws.Range("J91").Value = 1222.4286
$ws.Range("K91").Value = 2835
$ws.Range("L91").Value = 1222.4286
$ws.Range("M91").Value = -1431
$ws.Range("N91").Value = -4030.4286
$ws.Range("H135").Value = 4936
$ws.Range("I135").Value = 3081.45
$ws.Range("J135").Value = 10234.714
$ws.Range("K135").Value = 27733.05
$ws.Range("L135").Value = 92112.42600000001
$ws.Range("M135").Value = -25198.05
$ws.Range("N135").Value = -97182.42600000001
$ws.Range("H137").Value = 10527691
$ws.Range("I137").Value = 13334692
$ws.Range("J137").Value = 1437
$ws.Range("K137").Value = 40004076
$ws.Range("L137").Value = 4311
$ws.Range("M137").Value = -40001526
$ws.Range("N137").Value = -9411
$ws.Range("H138").Value = 7452.316
$ws.Range("I138").Value = 3000
$ws.Range("J138").Value = 7699.6665
$ws.Range("K138").Value = 9000
$ws.Range("L138").Value = 23098.9995
$ws.Range("M138").Value = -3860
$ws.Range("N138").Value = -33378.99950000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5109427.5
$ws.Range("I61").Value = 1784541.9
$ws.Range("K61").Value = 1784541.9
$ws.Range("M61").Value = -1784329.9
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 1362460.1
$ws.Range("I74").Value = 1613235.1
$ws.Range("K74").Value = 1613235.1
$ws.Range("M74").Value = -1612361.1
$ws.Range("H77").Value = 1362460.1
$ws.Range("I77").Value = 1613235.1
$ws.Range("K77").Value = 8066175.5
$ws.Range("M77").Value = -8061807.5
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H136").Value = 5109427.5
$ws.Range("I136").Value = 1784541.9
$ws.Range("K136").Value = 5353625.699999999
$ws.Range("M136").Value = -5351075.699999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H86").Value = 2989.7778
$ws.Range("I86").Value = 2565.4285
$ws.Range("J86").Value = 4475
$ws.Range("K86").Value = 2565.4285
$ws.Range("L86").Value = 4475
$ws.Range("M86").Value = -1442.4285
$ws.Range("N86").Value = -6721
$ws.Range("H89").Value = 2989.7778
$ws.Range("I89").Value = 2565.4285
$ws.Range("J89").Value = 4475
$ws.Range("K89").Value = 12827.1425
$ws.Range("L89").Value = 22375
$ws.Range("M89").Value = -7211.1425
$ws.Range("N89").Value = -33607
$ws.Range("H94").Value = 2145.1
$ws.Range("I94").Value = 1926.9131
$ws.Range("J94").Value = 2862
$ws.Range("K94").Value = 1926.9131
$ws.Range("L94").Value = 2862
$ws.Range("M94").Value = -1475.9131
$ws.Range("N94").Value = -3764
$ws.Range("H123").Value = 82500
$ws.Range("J123").Value = 82500
$ws.Range("L123").Value = 82500
$ws.Range("N123").Value = -92300

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1001999.5
$ws.Range("I6").Value = 1001999.5
$ws.Range("K6").Value = 1001999.5
$ws.Range("M6").Value = -1001886.5
$ws.Range("H58").Value = 10060286
$ws.Range("I58").Value = 16668445
$ws.Range("J58").Value = 3452126.5
$ws.Range("K58").Value = 16668445
$ws.Range("L58").Value = 3452126.5
$ws.Range("M58").Value = -16668242
$ws.Range("N58").Value = -3452532.5
$ws.Range("H99").Value = 13168.4
$ws.Range("I99").Value = 15618
$ws.Range("J99").Value = 5819.6
$ws.Range("K99").Value = 15618
$ws.Range("L99").Value = 5819.6
$ws.Range("M99").Value = -14120
$ws.Range("N99").Value = -8815.6
$ws.Range("H105").Value = 6434.625
$ws.Range("I105").Value = 4863.6
$ws.Range("J105").Value = 30000
$ws.Range("K105").Value = 4863.6
$ws.Range("L105").Value = 30000
$ws.Range("M105").Value = -3116.6
$ws.Range("N105").Value = -33494
$ws.Range("H107").Value = 727.75
$ws.Range("J107").Value = 597.8333
$ws.Range("L107").Value = 597.8333
$ws.Range("N107").Value = -4437.8333
$ws.Range("H126").Value = 13168.4
$ws.Range("I126").Value = 15618
$ws.Range("J126").Value = 5819.6
$ws.Range("K126").Value = 46854
$ws.Range("L126").Value = 17458.8
$ws.Range("M126").Value = -44384
$ws.Range("N126").Value = -22398.8
$ws.Range("H136").Value = 10060286
$ws.Range("I136").Value = 16668445
$ws.Range("J136").Value = 3452126.5
$ws.Range("K136").Value = 50005335
$ws.Range("L136").Value = 10356379.5
$ws.Range("M136").Value = -50002785
$ws.Range("N136").Value = -10361479.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4768.8
$ws.Range("J68").Value = 5040.902
$ws.Range("L68").Value = 15122.706
$ws.Range("N68").Value = -16744.706
$ws.Range("H71").Value = 4768.8
$ws.Range("J71").Value = 5040.902
$ws.Range("L71").Value = 45368.118
$ws.Range("N71").Value = -53480.118
$ws.Range("H98").Value = 696.1667
$ws.Range("I98").Value = 461
$ws.Range("K98").Value = 1383
$ws.Range("M98").Value = 115
$ws.Range("H128").Value = 281778.6
$ws.Range("I128").Value = 281778.6
$ws.Range("K128").Value = 845335.7999999999
$ws.Range("M128").Value = -840355.7999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2382.4783
$ws.Range("I126").Value = 2373.45
$ws.Range("K126").Value = 7120.349999999999
$ws.Range("M126").Value = -4650.349999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3728.625
$ws.Range("I22").Value = 2308.3333
$ws.Range("J22").Value = 4056.3845
$ws.Range("K22").Value = 2308.3333
$ws.Range("L22").Value = 4056.3845
$ws.Range("M22").Value = -2013.3333
$ws.Range("N22").Value = -4646.3845
$ws.Range("H27").Value = 3728.625
$ws.Range("I27").Value = 2308.3333
$ws.Range("J27").Value = 4056.3845
$ws.Range("K27").Value = 2308.3333
$ws.Range("L27").Value = 4056.3845
$ws.Range("M27").Value = -2201.3333
$ws.Range("N27").Value = -4270.3845
$ws.Range("H40").Value = 2899.889
$ws.Range("I40").Value = 2856.2856
$ws.Range("J40").Value = 3052.5
$ws.Range("K40").Value = 2856.2856
$ws.Range("L40").Value = 3052.5
$ws.Range("M40").Value = -2720.2856
$ws.Range("N40").Value = -3324.5
$ws.Range("H43").Value = 4235980
$ws.Range("I43").Value = 80000
$ws.Range("J43").Value = 5274975
$ws.Range("K43").Value = 80000
$ws.Range("L43").Value = 5274975
$ws.Range("M43").Value = -79807
$ws.Range("N43").Value = -5275361
$ws.Range("H82").Value = 4297.3335
$ws.Range("I82").Value = 4359.8
$ws.Range("J82").Value = 3985
$ws.Range("K82").Value = 4359.8
$ws.Range("L82").Value = 3985
$ws.Range("M82").Value = -3998.8
$ws.Range("N82").Value = -4707
$ws.Range("H85").Value = 4297.3335
$ws.Range("I85").Value = 4359.8
$ws.Range("J85").Value = 3985
$ws.Range("K85").Value = 4359.8
$ws.Range("L85").Value = 3985
$ws.Range("M85").Value = -3111.8
$ws.Range("N85").Value = -6481

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 491.58334
$ws.Range("I100").Value = 490.81818
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 981.63636
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -440.63636
$ws.Range("N100").Value = -2082
$ws.Range("H136").Value = 1787970.1
$ws.Range("I136").Value = 831999
$ws.Range("K136").Value = 2495997
$ws.Range("M136").Value = -2493447
